# daily auto push: 2025-10-09 02:00 UTC
# Append the new day's data row (row 83) to the bottom of the log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 83

# Column A holds a literal "yyyy/mm/dd" text label (same as every other
# row in the sheet), not a real date. Writing that string straight into a
# General-formatted cell would make Excel auto-convert it into a date
# serial number, so the cell is marked as Text first, then reset back to
# the default (General) look once the literal text is safely stored.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/09"
$dateCell.ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "木"
$ws.Cells.Item($newRow, 3).Value = 8
$ws.Cells.Item($newRow, 4).Value = 201
